$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each value below is set with a leading apostrophe so Excel stores it as
# literal Text (matching the source inlineStr cells) instead of re-parsing
# numeric-looking strings into Number cells (which would drop trailing zeros,
# e.g. "0.1200" -> 0.12).


# Row 2
$ws.Range("D2").Value = "'246.56"

# Row 3
$ws.Range("D3").Value = "'26.57"

# Row 4
$ws.Range("D4").Value = "'5.092"

# Row 5
$ws.Range("D5").Value = "'0.05614"

# Row 6
$ws.Range("D6").Value = "'6.478"

# Row 7
$ws.Range("D7").Value = "'0.8133"

# Row 8
$ws.Range("D8").Value = "'0.8443"

# Row 9
$ws.Range("B9").Value = "'BitrueCoin"
$ws.Range("C9").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D9").Value = "'0.02853"
$ws.Range("E9").Value = "'8BitrueCoinBTR"

# Row 10
$ws.Range("B10").Value = "'BitMartToken"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D10").Value = "'0.09388"
$ws.Range("E10").Value = "'9BitMartTokenBMX"

# Row 11
$ws.Range("B11").Value = "'BitForexToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D11").Value = "'0.001518"
$ws.Range("E11").Value = "'10BitForexTokenBF"

# Row 12
$ws.Range("B12").Value = "'One"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D12").Value = "'0.0005962"
$ws.Range("E12").Value = "'11OneONE"

# Row 13
$ws.Range("D13").Value = "'0.006109"

# Row 14
$ws.Range("D14").Value = "'3.590"

# Row 15
$ws.Range("D15").Value = "'3.010"

# Row 16
$ws.Range("D16").Value = "'2.055"

# Row 18
$ws.Range("D18").Value = "'0.1339"

# Row 19
$ws.Range("D19").Value = "'0.06969"

# Row 20
$ws.Range("D20").Value = "'0.03143"

# Row 21
$ws.Range("D21").Value = "'0.1320"

# Row 22
$ws.Range("D22").Value = "'3.751"

# Row 23
$ws.Range("D23").Value = "'0.04655"

# Row 25
$ws.Range("D25").Value = "'0.001252"

# Row 27
$ws.Range("D27").Value = "'0.00009604"

# Row 28
$ws.Range("E28").Value = "'27UpBotsUBXTBestin24h"

# Row 40
$ws.Range("D40").Value = "'0.03669"

# Row 41
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1358"
$ws.Range("E41").Value = "'40BKEXTokenBKK"

# Row 42
$ws.Range("B42").Value = "'KickToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006183"
$ws.Range("E42").Value = "'41KickTokenKICK"

# Row 43
$ws.Range("D43").Value = "'0.002661"

# Row 44
$ws.Range("D44").Value = "'0.008904"

# Row 47
$ws.Range("D47").Value = "'0.1200"

# Row 48
$ws.Range("D48").Value = "'0.002516"
